$wb = $excel.ActiveWorkbook

# Reference sheet that already has the desired header style (bold, centered,
# bordered) so the new sheets' header rows match the look of the existing
# ones.
$headerSource = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------
# New sheet: "ODI Batting Extra"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

# Copy header formatting (bold font, centered alignment, thin border) from
# an existing header row, then overwrite the header text.
$headerSource.Range("A1:F1").Copy()
$battingExtra.Range("A1:F1").PasteSpecial(-4122)

$battingExtra.Range("A1").Value = "MATCH_CODE"
$battingExtra.Range("B1").Value = "BATTING_POSITION"
$battingExtra.Range("C1").Value = "NUM_4"
$battingExtra.Range("D1").Value = "NUM_6"
$battingExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$battingExtra.Range("F1").Value = "MAN_OF_MATCH"

# Data row. Most values are stored as plain text (matching the rest of the
# workbook), but BATTING_POSITION is a genuine number.
$battingExtra.Range("A2").NumberFormat = "@"
$battingExtra.Range("A2").Value = "3996"

$battingExtra.Range("B2").Value = 6

$battingExtra.Range("C2:F2").NumberFormat = "@"
$battingExtra.Range("C2").Value = "0"
$battingExtra.Range("D2").Value = "0"
$battingExtra.Range("E2").Value = "8.11%"
$battingExtra.Range("F2").Value = "NO"

# ---------------------------------------------------------------------
# New sheet: "ODI Bowling Extra"
# ---------------------------------------------------------------------
$bowlingExtra = $wb.Worksheets.Add($null, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"

$headerSource.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$bowlingExtra.Range("A2:B2").NumberFormat = "@"
$bowlingExtra.Range("A2").Value = "3996"
$bowlingExtra.Range("B2").Value = "0"
# PERCENT_WICKETS_OF_ALL is blank for this row.

# Restore the originally active sheet/selection so the workbook-level view
# state is unchanged by adding the new sheets.
$null = $wb.Worksheets.Item(1).Activate()
